$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "vaishali.kh2310@gmail.com"
$ws.Range("B21").Value = "Login"
$ws.Range("C21").Value = "2025-06-15 22:17:14"
